# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q3" and "总计",
#    populated with the Q1-2022 fund-holding data.
# 2. Update the "总计" (totals) sheet with a new row summarising the
#    2022-Q1 sheet, keeping the existing 2021-Q3 row below it.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)

# --- 1. Create the new "2022-Q1" sheet right after "2021-Q3" ---------------
$newSheet = $wb.Worksheets.Add($null, $sheet1, 1)
$newSheet.Name = "2022-Q1"

# Worksheet references captured before the sheet count changed can become
# stale, so re-resolve "总计" by name now that the new sheet exists.
$total = $wb.Worksheets.Item("总计")

# Borrow the existing formatting (font/border/alignment) used on "总计" for
# the header row and the first (index) column, so the new sheet matches the
# established look instead of Excel's bare defaults.
$total.Range("B1").Copy($newSheet.Range("B1:H1"))
$total.Range("A2").Copy($newSheet.Range("A2:A3"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'920003"
$newSheet.Range("C2").Value = "中金新锐股票A"
$newSheet.Range("D2").Value = "'24.64"
$newSheet.Range("E2").Value = "'92.76"
$newSheet.Range("F2").Value = "'3.80"
$newSheet.Range("G2").Value = "'0.9363"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'920923"
$newSheet.Range("C3").Value = "中金新锐股票C"
$newSheet.Range("D3").Value = "'3.94"
$newSheet.Range("E3").Value = "'92.76"
$newSheet.Range("F3").Value = "'3.80"
$newSheet.Range("G3").Value = "'0.1497"
$newSheet.Range("H3").Value = 4

# --- 2. Add the 2022-Q1 summary row to "总计", above the 2021-Q3 row -------
$total.Rows(2).Insert()

# The freshly inserted row picks up a default/blank style; reset it and pull
# in the same index-column formatting used by the row it displaced.
$total.Range("A2:D2").ClearFormats()
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 1.09
